{"js": "// Taller 1.docx \u2014 update the placeholder file-name the students must use\n// when renaming their PDF: \"NombreApellidoEVIMP2020II\" -> \"NombreApellidoEEC2021I\"\n// (i.e. \"VIMP2020II\" is replaced by \"EC2021I\" inside the quoted instructions).\n\nconst body = context.document.body;\n\n// Locate the exact (unique) run of text that holds the old placeholder.\nconst results = body.search(\"NombreApellidoEVIMP2020II\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Expected text \"NombreApellidoEVIMP2020II\" not found.');\n}\n\n// Replace it in place with the new placeholder text.\nresults.items[0].insertText(\"NombreApellidoEEC2021I\", \"Replace\");\nawait context.sync();\n", "ps1": "# Taller 1.docx \u2014 update the placeholder file-name the students must use\n# when renaming their PDF: \"NombreApellidoEVIMP2020II\" -> \"NombreApellidoEEC2021I\"\n# (i.e. \"VIMP2020II\" is replaced by \"EC2021I\" inside the quoted instructions).\n\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"NombreApellidoEVIMP2020II\"\n$find.MatchCase = $true\n$find.Execute() | Out-Null\n\nif ($find.Found) {\n    $range.Text = \"NombreApellidoEEC2021I\"\n}\n"}
